$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Cells.Item(8, 8).Value = 15151754
$ws.Cells.Item(8, 9).Value = 33333528
$ws.Cells.Item(8, 10).Value = 276.08334
$ws.Cells.Item(8, 11).Value = 100000584
$ws.Cells.Item(8, 12).Value = 828.2500200000001
$ws.Cells.Item(8, 13).Value = -100000445
$ws.Cells.Item(8, 14).Value = -1106.25002
# Row 43
$ws.Cells.Item(43, 8).Value = 7935.8335
$ws.Cells.Item(43, 9).Value = 4402
$ws.Cells.Item(43, 11).Value = 4402
$ws.Cells.Item(43, 13).Value = -4333
# Row 55
$ws.Cells.Item(55, 8).Value = 120.1
$ws.Cells.Item(55, 9).Value = 168.4
$ws.Cells.Item(55, 10).Value = 71.8
$ws.Cells.Item(55, 11).Value = 168.4
$ws.Cells.Item(55, 12).Value = 71.8
$ws.Cells.Item(55, 13).Value = 45.59999999999999
$ws.Cells.Item(55, 14).Value = -499.8
# Row 61
$ws.Cells.Item(61, 8).Value = 1109.875
$ws.Cells.Item(61, 9).Value = 911.2857
$ws.Cells.Item(61, 11).Value = 2733.8571
$ws.Cells.Item(61, 13).Value = -2561.8571
# Row 135
$ws.Cells.Item(135, 8).Value = 3955.1333
$ws.Cells.Item(135, 9).Value = 3973.2693
$ws.Cells.Item(135, 11).Value = 35759.4237
$ws.Cells.Item(135, 13).Value = -33224.4237
# Row 138
$ws.Cells.Item(138, 8).Value = 3786.53
$ws.Cells.Item(138, 10).Value = 5116.82
$ws.Cells.Item(138, 12).Value = 15350.46
$ws.Cells.Item(138, 14).Value = -25630.46
# Row 141
$ws.Cells.Item(141, 8).Value = 8494.6
$ws.Cells.Item(141, 9).Value = 8699.625
$ws.Cells.Item(141, 11).Value = 26098.875
$ws.Cells.Item(141, 13).Value = -20918.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5522.6196
$ws.Cells.Item(32, 9).Value = 5489.789
$ws.Cells.Item(32, 11).Value = 5489.789
$ws.Cells.Item(32, 13).Value = -5202.789
# Row 45
$ws.Cells.Item(45, 8).Value = 76010.03999999999
$ws.Cells.Item(45, 9).Value = 94666.67999999999
$ws.Cells.Item(45, 11).Value = 94666.67999999999
$ws.Cells.Item(45, 13).Value = -94289.67999999999
# Row 61
$ws.Cells.Item(61, 8).Value = 5177.822
$ws.Cells.Item(61, 9).Value = 5230.4
$ws.Cells.Item(61, 11).Value = 5230.4
$ws.Cells.Item(61, 13).Value = -5018.4
# Row 122
$ws.Cells.Item(122, 8).Value = 942050.1
$ws.Cells.Item(122, 9).Value = 3985.3333
$ws.Cells.Item(122, 11).Value = 11955.9999
$ws.Cells.Item(122, 13).Value = -9505.999899999999
# Row 136
$ws.Cells.Item(136, 8).Value = 5177.822
$ws.Cells.Item(136, 9).Value = 5230.4
$ws.Cells.Item(136, 11).Value = 15691.2
$ws.Cells.Item(136, 13).Value = -13141.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 1131.8667
$ws.Cells.Item(107, 9).Value = 998.1818
$ws.Cells.Item(107, 11).Value = 998.1818
$ws.Cells.Item(107, 13).Value = 921.8182
# Row 134
$ws.Cells.Item(134, 8).Value = 1466.5122
$ws.Cells.Item(134, 10).Value = 3453.1667
$ws.Cells.Item(134, 12).Value = 10359.5001
$ws.Cells.Item(134, 14).Value = -15429.5001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 7150.1514
$ws.Cells.Item(31, 9).Value = 7231.8335
$ws.Cells.Item(31, 10).Value = 6333.3335
$ws.Cells.Item(31, 11).Value = 7231.8335
$ws.Cells.Item(31, 12).Value = 6333.3335
$ws.Cells.Item(31, 13).Value = -6936.8335
$ws.Cells.Item(31, 14).Value = -6923.3335
# Row 34
$ws.Cells.Item(34, 8).Value = 7150.1514
$ws.Cells.Item(34, 9).Value = 7231.8335
$ws.Cells.Item(34, 10).Value = 6333.3335
$ws.Cells.Item(34, 11).Value = 7231.8335
$ws.Cells.Item(34, 12).Value = 6333.3335
$ws.Cells.Item(34, 13).Value = -7029.8335
$ws.Cells.Item(34, 14).Value = -6737.3335
# Row 58
$ws.Cells.Item(58, 8).Value = 2453.0527
$ws.Cells.Item(58, 9).Value = 1729
$ws.Cells.Item(58, 11).Value = 1729
$ws.Cells.Item(58, 13).Value = -1526
# Row 132
$ws.Cells.Item(132, 8).Value = 13467.889
$ws.Cells.Item(132, 9).Value = 18202
$ws.Cells.Item(132, 11).Value = 54606
$ws.Cells.Item(132, 13).Value = -52076
# Row 136
$ws.Cells.Item(136, 8).Value = 2453.0527
$ws.Cells.Item(136, 9).Value = 1729
$ws.Cells.Item(136, 11).Value = 5187
$ws.Cells.Item(136, 13).Value = -2637
# Row 141
$ws.Cells.Item(141, 8).Value = 129717.11
$ws.Cells.Item(141, 10).Value = 134389.23
$ws.Cells.Item(141, 12).Value = 134389.23
$ws.Cells.Item(141, 14).Value = -144749.23

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 43.75
$ws.Cells.Item(2, 9).Value = 43.3
$ws.Cells.Item(2, 11).Value = 259.8
$ws.Cells.Item(2, 13).Value = -146.8
# Row 131
$ws.Cells.Item(131, 8).Value = 15875677
$ws.Cells.Item(131, 10).Value = 1781.9348
$ws.Cells.Item(131, 12).Value = 5345.8044
$ws.Cells.Item(131, 14).Value = -15425.8044
# Row 139
$ws.Cells.Item(139, 8).Value = 1156392.9
$ws.Cells.Item(139, 9).Value = 1766425.9
$ws.Cells.Item(139, 10).Value = 4108.222
$ws.Cells.Item(139, 11).Value = 5299277.699999999
$ws.Cells.Item(139, 12).Value = 12324.666
$ws.Cells.Item(139, 13).Value = -5294137.699999999
$ws.Cells.Item(139, 14).Value = -22604.666
# Row 140
$ws.Cells.Item(140, 8).Value = 3049.2
$ws.Cells.Item(140, 9).Value = 2967.9167
$ws.Cells.Item(140, 11).Value = 8903.750100000001
$ws.Cells.Item(140, 13).Value = -3723.750100000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Cells.Item(113, 8).Value = 12276.091
$ws.Cells.Item(113, 10).Value = 3325.5
$ws.Cells.Item(113, 12).Value = 3325.5
$ws.Cells.Item(113, 14).Value = -7665.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 18198.842
$ws.Cells.Item(40, 9).Value = 19191.857
$ws.Cells.Item(40, 11).Value = 19191.857
$ws.Cells.Item(40, 13).Value = -19055.857
# Row 68
$ws.Cells.Item(68, 8).Value = 6417.0835
$ws.Cells.Item(68, 10).Value = 9007.429
$ws.Cells.Item(68, 12).Value = 9007.429
$ws.Cells.Item(68, 14).Value = -10505.429
# Row 71
$ws.Cells.Item(71, 8).Value = 6417.0835
$ws.Cells.Item(71, 10).Value = 9007.429
$ws.Cells.Item(71, 12).Value = 45037.145
$ws.Cells.Item(71, 14).Value = -52525.145
# Row 87
$ws.Cells.Item(87, 8).Value = 55200
$ws.Cells.Item(87, 10).Value = 55200
$ws.Cells.Item(87, 12).Value = 55200
$ws.Cells.Item(87, 14).Value = -57446
# Row 90
$ws.Cells.Item(90, 8).Value = 55200
$ws.Cells.Item(90, 10).Value = 55200
$ws.Cells.Item(90, 12).Value = 165600
$ws.Cells.Item(90, 14).Value = -176832
# Row 136
$ws.Cells.Item(136, 8).Value = 4089.8667
$ws.Cells.Item(136, 9).Value = 2949.8684
$ws.Cells.Item(136, 10).Value = 6058.9546
$ws.Cells.Item(136, 11).Value = 8849.6052
$ws.Cells.Item(136, 12).Value = 18176.8638
$ws.Cells.Item(136, 13).Value = -6299.6052
$ws.Cells.Item(136, 14).Value = -23276.8638

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 229867.17
$ws.Cells.Item(62, 9).Value = 425563.44
$ws.Cells.Item(62, 11).Value = 425563.44
$ws.Cells.Item(62, 13).Value = -424939.44
# Row 65
$ws.Cells.Item(65, 8).Value = 229867.17
$ws.Cells.Item(65, 9).Value = 425563.44
$ws.Cells.Item(65, 11).Value = 2127817.2
$ws.Cells.Item(65, 13).Value = -2124697.2
# Row 122
$ws.Cells.Item(122, 8).Value = 5584.5293
$ws.Cells.Item(122, 9).Value = 3625.4062
$ws.Cells.Item(122, 11).Value = 10876.2186
$ws.Cells.Item(122, 13).Value = -8426.2186
# Row 126
$ws.Cells.Item(126, 8).Value = 23232.61
$ws.Cells.Item(126, 10).Value = 5095.5557
$ws.Cells.Item(126, 12).Value = 15286.6671
$ws.Cells.Item(126, 14).Value = -20226.6671
# Row 132
$ws.Cells.Item(132, 8).Value = 8813.456
$ws.Cells.Item(132, 9).Value = 9600.5
$ws.Cells.Item(132, 11).Value = 28801.5
$ws.Cells.Item(132, 13).Value = -26271.5
